$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 with the new control-point values
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 6

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 6

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 3

# Remove the now-obsolete rows 5-9 entirely
$ws.Range("A5:B9").EntireRow.Delete()
